$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# fix orga: Organisation value "CMML" -> "HHU CMML"
$ws.Range("B6").Value = "HHU CMML"

# add author: Author Last Name / Author First Name
$ws.Range("B17").Value = "Brilhaus"
$ws.Range("B18").Value = "Dominik"
